$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 152 (pushes the existing 152..216 down to 153..217)
$ws.Rows.Item(152).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A152").Value = 4
$ws.Range("B152").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C152").Value = "Los Lagos"
$ws.Range("D152").Value = 44704
$ws.Range("E152").Value = 10
$ws.Range("F152").Value = "Fruta"
$ws.Range("G152").Value = 100103
$ws.Range("H152").Value = "Frutos de hueso (carozo)"
$ws.Range("I152").Value = 100103002
$ws.Range("J152").Value = "Ciruela"
$ws.Range("K152").Value = "Black Amber"
$ws.Range("L152").Value = "Segunda"
$ws.Range("M152").Value = 300
$ws.Range("N152").Value = 13000
$ws.Range("O152").Value = 13000
$ws.Range("P152").Value = 13000
$ws.Range("Q152").Value = "$/caja 15 kilos granel"
$ws.Range("R152").Value = "Provincia de Curicó"
$ws.Range("S152").Value = 867
$ws.Range("T152").Value = 15
